# Insert two new "Title and Content" slides right after the title slide
# (slide 1), pushing the existing "Dungeon Jump" / "Overworld" /
# "Dungeon Instances" slides down from positions 2-4 to positions 4-6.

$p = $ppt.ActivePresentation
$layout = $p.SlideMaster.CustomLayouts.Item(2)   # "Title and Content"

# --- New slide 2: "Git directory structure" ---------------------------
$gitDirSlide = $p.Slides.AddSlide(2, $layout)
$gitDirSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Git directory structure"
$gitDirSlide.Shapes.Item(2).TextFrame.TextRange.Text = "To be demonstrated"

# --- New slide 3: "Team member Git contributions" ----------------------
$contribSlide = $p.Slides.AddSlide(3, $layout)
$contribSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Team member Git contributions"

$body = $contribSlide.Shapes.Item(2).TextFrame.TextRange
$body.Text = "To be demonstrated`r/src readme`r/doc Git manual`r/doc Champion documents"
$body.Paragraphs(2).IndentLevel = 2
$body.Paragraphs(3).IndentLevel = 2
$body.Paragraphs(4).IndentLevel = 2
